# Fruta / hortaliza, semanal
# A new weekly price-report row is prepended to the data table (row 7),
# pushing the existing rows 7-34 down to 8-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 7 - this shifts rows 7..34 down to 8..35
# and extends the sheet dimension from A1:R34 to A1:R35.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with this week's record.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44819
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112045
$ws.Range("G7").Value = "Zapallo"
$ws.Range("H7").Value = "Camote"
$ws.Range("I7").Value = "1a nueva(o)"
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 950
$ws.Range("M7").Value = 925
$ws.Range("N7").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O7").Value = "Perú"
$ws.Range("P7").Value = 925
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
